$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 6: Inscritos (E6) 54 -> 56
$ws.Range("E6").Value = 56

# Row 7: Inscritos (E7) 30 -> 31, Pagos (F7) 17 -> 19, Inscrições homologadas (H7) 18 -> 20
$ws.Range("E7").Value = 31
$ws.Range("F7").Value = 19
$ws.Range("H7").Value = 20

# Row 12: Inscritos (E12) 29 -> 31, Pagos (F12) 12 -> 13, Inscrições homologadas (H12) 14 -> 15
$ws.Range("E12").Value = 31
$ws.Range("F12").Value = 13
$ws.Range("H12").Value = 15

# Row 16: Inscritos (E16) 310 -> 312, Pagos (F16) 85 -> 88, Inscrições homologadas (H16) 172 -> 175
$ws.Range("E16").Value = 312
$ws.Range("F16").Value = 88
$ws.Range("H16").Value = 175

# Row 18: Pagos (F18) 28 -> 30, Inscrições homologadas (H18) 51 -> 53
$ws.Range("F18").Value = 30
$ws.Range("H18").Value = 53
